$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update cell text content (drives sharedStrings regeneration) ---
$ws.Range("F5").Value2 = "O Plano de Projeto não está no repositório para verificação de planejamento de baselines. Não foram criadas baselines."
$ws.Range("F9").Value2 = "Ainda não foram feitas auditorias, pois ainda não existem baselines."
$ws.Range("C10").Value2 = "NA"
$ws.Range("F10").Value2 = "Não existem documentos do processo GCO no repositório além do Plano de Configuração. O Plano de Projeto não está disponível no repositório para verificação dos artefatos que devem ser entregues pela GCO e o Plano de Configuração não contempla essa informação."
$ws.Range("F13").Value2 = "O Plano de Projeto não está disponível no Repositório. O Plano de Configuração não foi aprovado."
$ws.Range("C14").Value2 = "Não"
$ws.Range("F14").Value2 = "O Plano de Configuração não foi aprovado."

# --- Column widths ---
$ws.Columns.Item(2).ColumnWidth = 60.85546875
$ws.Columns.Item(4).ColumnWidth = 19.5703125
$ws.Columns.Item(5).ColumnWidth = 14.85546875
$ws.Columns.Item(6).ColumnWidth = 50.42578125

# --- Row heights ---
$ws.Rows.Item(2).RowHeight = 54
$ws.Rows.Item(3).RowHeight = 48.75
$ws.Rows.Item(4).RowHeight = 48
$ws.Rows.Item(5).RowHeight = 50.25
$ws.Rows.Item(6).RowHeight = 48
$ws.Rows.Item(7).RowHeight = 44.25
$ws.Rows.Item(8).RowHeight = 42.75
$ws.Rows.Item(9).RowHeight = 49.5
$ws.Rows.Item(10).RowHeight = 72
$ws.Rows.Item(11).RowHeight = 39
$ws.Rows.Item(12).RowHeight = 39.75
$ws.Rows.Item(13).RowHeight = 42.75
$ws.Rows.Item(14).RowHeight = 42.75

# --- Wrap text + vertical center for header row and body rows ---
$ws.Range("A1:F14").WrapText = $true
$ws.Range("A1:F14").VerticalAlignment = -4108

# --- Selection ---
$ws.Range("G4").Select()
